$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.84"
$ws.Range("E2").Value = "'4.85%"
$ws.Range("D3").Value = "'36.14"
$ws.Range("E3").Value = "'16.73%"
$ws.Range("D4").Value = "'5.063"
$ws.Range("E4").Value = "'2.43%"
$ws.Range("D5").Value = "'0.07849"
$ws.Range("E5").Value = "'6.25%"
$ws.Range("D6").Value = "'2.290"
$ws.Range("E6").Value = "'1.48%"
$ws.Range("D7").Value = "'8.069"
$ws.Range("E7").Value = "'4.68%"
$ws.Range("D8").Value = "'4.000"
$ws.Range("E8").Value = "'6.60%"
$ws.Range("D9").Value = "'0.9270"
$ws.Range("E9").Value = "'0.93%"
$ws.Range("D10").Value = "'0.1012"
$ws.Range("E10").Value = "'10.15%"
$ws.Range("E11").Value = "'7.24%"
$ws.Range("D12").Value = "'0.08590"
$ws.Range("E12").Value = "'3.74%"
$ws.Range("D13").Value = "'0.03406"
$ws.Range("E13").Value = "'5.74%"
$ws.Range("D14").Value = "'0.09905"
$ws.Range("E14").Value = "'-0.87%"
$ws.Range("D15").Value = "'0.001478"
$ws.Range("E15").Value = "'-1.49%"
$ws.Range("D16").Value = "'0.04674"
$ws.Range("E16").Value = "'3.59%"
$ws.Range("D17").Value = "'0.005610"
$ws.Range("E17").Value = "'-3.51%"
$ws.Range("E18").Value = "'0.42%"
$ws.Range("D19").Value = "'2.097"
$ws.Range("E19").Value = "'-0.30%"
$ws.Range("D20").Value = "'0.3433"
$ws.Range("E20").Value = "'3.11%"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("E21").Value = "'3.19%"
$ws.Range("D22").Value = "'4.547"
$ws.Range("E22").Value = "'9.73%"
$ws.Range("D23").Value = "'0.2224"
$ws.Range("E23").Value = "'4.84%"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'2.09%"
$ws.Range("D25").Value = "'0.004482"
$ws.Range("E25").Value = "'6.41%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'0.10%"
$ws.Range("D27").Value = "'0.0003003"
$ws.Range("E27").Value = "'-11.43%"
$ws.Range("D39").Value = "'0.01750"
$ws.Range("E39").Value = "'10.55%"
$ws.Range("D40").Value = "'0.04694"
$ws.Range("E40").Value = "'3.25%"
$ws.Range("D41").Value = "'0.007900"
$ws.Range("E41").Value = "'6.64%"
$ws.Range("D42").Value = "'0.1413"
$ws.Range("E42").Value = "'5.48%"
$ws.Range("D43").Value = "'0.008814"
$ws.Range("E43").Value = "'-10.51%"
$ws.Range("D44").Value = "'0.002211"
$ws.Range("E44").Value = "'-0.29%"
$ws.Range("D45").Value = "'0.009182"
$ws.Range("E45").Value = "'7.66%"
$ws.Range("D46").Value = "'0.00005980"
$ws.Range("E46").Value = "'-2.01%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.22%"
$ws.Range("E48").Value = "'132.51%"
$ws.Range("D49").Value = "'0.002692"
$ws.Range("E49").Value = "'34.78%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.22%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.22%"
